$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete years 2004-2009 (old rows 2-7). This shifts the
# remaining rows (old 2010-2020 data) up to become new rows 2-12.
$ws.Range("A2:A7").EntireRow.Delete()

# Copy formatting (style) of the last existing data row (2020, now row 12)
# down into the new row 13 before filling in the 2021 figures.
$ws.Range("A12").Copy($ws.Range("A13"))

# Populate the new row for 2021.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 274.7809924434
$ws.Range("C13").Value = 988.3782761433
$ws.Range("E13").Value = 537.5976949039
